$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8815737962722778
$ws.Range("B1").Value = 0.7916948199272156
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.53313934803009
$ws.Range("E1").Value = 0.9393343925476074
